$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update rows affected by the daily data refresh / re-sort described in the
# commit "Update countries & provincias Spain". Nepal and Libia move up in the
# ranking (pushing Suecia/Guatemala/Chequia and Kenia/Afganistan/Irlanda down a
# row respectively), and Islas Malvinas / Montserrat swap places for a tie.

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 7777014
$ws.Range("C4").Value = 790
$ws.Range("D4").Value = 4984464
$ws.Range("E4").Value = 2575732
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 34
$ws.Range("H4").Value = 216818

# Row 16: Iran
$ws.Range("A16").Value = "Iran"
$ws.Range("B16").Value = 488236
$ws.Range("C16").Value = 4392
$ws.Range("D16").Value = 399300
$ws.Range("E16").Value = 61048
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 230
$ws.Range("H16").Value = 27888

# Row 44: Emiratos Arabes Unidos
$ws.Range("A44").Value = "Emiratos Arabes Unidos"
$ws.Range("B44").Value = 102929
$ws.Range("C44").Value = 1089
$ws.Range("D44").Value = 93479
$ws.Range("E44").Value = 9012
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 438

# Row 45: Nepal
$ws.Range("A45").Value = "Nepal"
$ws.Range("B45").Value = 98617
$ws.Range("C45").Value = 4364
$ws.Range("D45").Value = 71343
$ws.Range("E45").Value = 26684
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 12
$ws.Range("H45").Value = 590

# Row 46: Suecia
$ws.Range("A46").Value = "Suecia"
$ws.Range("B46").Value = 96677
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 5892

# Row 47: Guatemala
$ws.Range("A47").Value = "Guatemala"
$ws.Range("B47").Value = 95704
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 84036
$ws.Range("E47").Value = 8333
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 3335

# Row 48: Chequia
$ws.Range("A48").Value = "Chequia"
$ws.Range("B48").Value = 95360
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 50767
$ws.Range("E48").Value = 43764
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 829

# Row 71: Azerbaiyan
$ws.Range("A71").Value = "Azerbaiyan"
$ws.Range("B71").Value = 41304
$ws.Range("C71").Value = 191
$ws.Range("D71").Value = 39012
$ws.Range("E71").Value = 1689
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 603

# Row 72: Libia
$ws.Range("A72").Value = "Libia"
$ws.Range("B72").Value = 40292
$ws.Range("C72").Value = 779
$ws.Range("D72").Value = 23130
$ws.Range("E72").Value = 16546
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 8
$ws.Range("H72").Value = 616

# Row 73: Kenia
$ws.Range("A73").Value = "Kenia"
$ws.Range("B73").Value = 39907
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 31659
$ws.Range("E73").Value = 7500
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 748

# Row 74: Afganistan
$ws.Range("A74").Value = "Afganistan"
$ws.Range("B74").Value = 39616
$ws.Range("C74").Value = 68
$ws.Range("D74").Value = 33058
$ws.Range("E74").Value = 5088
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 1470

# Row 75: Irlanda
$ws.Range("A75").Value = "Irlanda"
$ws.Range("B75").Value = 39584
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 23364
$ws.Range("E75").Value = 14404
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 1816

# Row 110: Uganda
$ws.Range("A110").Value = "Uganda"
$ws.Range("B110").Value = 9442
$ws.Range("C110").Value = 182
$ws.Range("D110").Value = 5781
$ws.Range("E110").Value = 3576
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 85

# Row 215: Islas Malvinas
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0

# Row 216: Montserrat
$ws.Range("A216").Value = "Montserrat"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 12
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 1
